# Development of the xls file with description of the DNS geometry.
#
# This script reproduces, via Excel COM automation, the edits that were made
# to dns_draft.xlsx: a "Parameter / Value / Units" header row was added to
# both the PSD_8Packs and PA sheets, the separate tube_width/tube_thickness
# rows were merged into a single tube_diameter row (with a value), a new
# angular_step parameter row was introduced on both sheets, the air_gap_width
# row on PSD_8Packs got a value + unit, and the window/selection state was
# refreshed.

$wb = $excel.ActiveWorkbook

$wsGeometry = $wb.Worksheets.Item("Geometry")
$wsPsd      = $wb.Worksheets.Item("PSD_8Packs")
$wsPa       = $wb.Worksheets.Item("PA")

# ---------------------------------------------------------------------
# Sheet "PSD_8Packs": rows, before any insert, are
#   1 pixels_per_tube  960
#   2 tubes_per_bank   8
#   3 air_gap_width
#   4 tube_length      1     m
#   5 tube_width
#   6 tube_thickness
#   7 tube_pressure    16    bar
#   8 tube_temperature 300   K
# ---------------------------------------------------------------------

# Turn the old "tube_width" row into the new "tube_diameter" row (value in
# metres) before the header row shifts everything down by one.
$wsPsd.Range("A5").Value = "tube_diameter"
$wsPsd.Range("B5").Value = 0.0127
$wsPsd.Range("C5").Value = "m"

# Insert the new header row at the top.
$wsPsd.Rows.Item(1).Insert()

$wsPsd.Range("A1").Value = "Parameter"
$wsPsd.Range("B1").Value = "Value"
$wsPsd.Range("C1").Value = "Units"
$wsPsd.Range("A1:C1").Font.Bold = $true

# air_gap_width (now row 4) gets a value + unit.
$wsPsd.Range("B4").Value = 1
$wsPsd.Range("C4").Value = "mm"

# "tube_thickness" (now row 7) becomes the new "angular_step" row.
$wsPsd.Range("A7").Value = "angular_step"
$wsPsd.Range("B7").Value = 1
$wsPsd.Range("C7").Value = "degree"

# ---------------------------------------------------------------------
# Sheet "PA": rows, before any insert, are
#   1 pixels_per_tube  1
#   2 tubes_per_bank   24
#   3 air_gap_width
#   4 tube_length      0.15  m
#   5 tube_width
#   6 tube_thickness
#   7 tube_pressure    16    bar
#   8 tube_temperature 300   K
# ---------------------------------------------------------------------

# Turn the old "tube_width" row into the new "tube_diameter" row.
$wsPa.Range("A5").Value = "tube_diameter"
$wsPa.Range("B5").Value = 0.0254
$wsPa.Range("C5").Value = "m"

# Insert the new header row at the top.
$wsPa.Rows.Item(1).Insert()

$wsPa.Range("A1").Value = "Parameter"
$wsPa.Range("B1").Value = "Value"
$wsPa.Range("C1").Value = "Units"
$wsPa.Range("A1:C1").Font.Bold = $true

# The old empty "air_gap_width" row (now row 4) is replaced outright by the
# new "angular_step" row -- PA has no air_gap_width parameter.
$wsPa.Range("A4").Value = "angular_step"
$wsPa.Range("B4").Value = 5
$wsPa.Range("C4").Value = "degree"

# "tube_thickness" (now row 7, empty) is no longer part of the table.
$wsPa.Rows.Item(7).Delete() | Out-Null

# ---------------------------------------------------------------------
# Window / selection state.
# ---------------------------------------------------------------------

# Reposition the workbook window.
$excel.ActiveWindow.Left = 3300
$excel.ActiveWindow.Top = 3220

# Selections on each sheet (also updates the active-sheet/tabSelected flag,
# so these run in an order that leaves "Geometry" as the active tab last).
$wsPsd.Range("A1:C1").Select() | Out-Null
$wsPa.Range("C5").Select() | Out-Null
$wsGeometry.Activate()
$wsGeometry.Range("H2").Select() | Out-Null
